$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value reads as a plain number (e.g. "309.38").
# Force them to Text format first so Excel stores the literal string instead of
# silently converting it to a numeric value (matches how the source data is stored).
$textCells = $ws.Range('D4,D5,D6,D7,D8,D9,D10,D11,D12,D14,D15,D16,D17,D18,D19,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51')
$textCells.NumberFormat = "@"

# --- Row-by-row Price / Volume(1h) refresh ---
$ws.Range('D2').Value = '26.889.10'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '1.843.45'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '309.38'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').Value = '0.4719'
$ws.Range('E7').Value = '  +0.88%  '
$ws.Range('D8').Value = '0.3678'
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('D9').Value = '0.07223'
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('D10').Value = '0.9233'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('D11').Value = '19.61'
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').Value = '0.07615'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '1.881.10'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('D14').Value = '5.307'
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').Value = '6.392'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '88.33'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '0.000008656'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '26.935.30'
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '14.56'
$ws.Range('E21').Value = '  +2.61%  '
$ws.Range('D22').Value = '5.033'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').Value = '10.66'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '1.914'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').Value = '152.11'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').Value = '18.15'
$ws.Range('E26').Value = '  +1.42%  '
$ws.Range('D27').Value = '2.000'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('D28').Value = '114.19'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').Value = '4.966'
$ws.Range('E29').Value = '  +3.38%  '
$ws.Range('D30').Value = '0.08831'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').Value = '3.280'
$ws.Range('E31').Value = '  +4.23%  '
$ws.Range('D32').Value = '0.7469'
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('D33').Value = '1.166'
$ws.Range('E33').Value = '  +3.84%  '
$ws.Range('D36').Value = '1.088'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '0.01948'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('D38').Value = '0.05253'
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').Value = '2.966'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('D40').Value = '0.5202'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').Value = '6.904'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('D42').Value = '0.1512'
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('D43').Value = '8.199'
$ws.Range('E43').Value = '  +2.58%  '
$ws.Range('D44').Value = '10.55'
$ws.Range('E44').Value = '  +5.90%  '
$ws.Range('D45').Value = '0.4693'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').Value = '1.006'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('D47').Value = '101.92'
$ws.Range('E47').Value = '  +2.98%  '
$ws.Range('D48').Value = '1.599'
$ws.Range('D49').Value = '65.46'
$ws.Range('E49').Value = '  +2.64%  '
$ws.Range('D50').Value = '0.06030'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '0.8831'
$ws.Range('E51').Value = '  +3.81%  '

# --- Rows 34/35: RenderToken and Filecoin swapped position in the ranking ---
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.492'
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '2.759'
$ws.Range('E35').Value = '  -0.27%  '
